# Generate Report for Archive
#
# The nightly localization-status report changed in two related ways:
#   1. Items that used to show the text "Ready for handoff" now show
#      "In Translation" (this text lives in the Overview sheet's zh-cn /
#      de-de columns, and in the "Status" column of each per-locale
#      sheet).
#   2. Because the new status string is shorter, the "Status" style
#      column(s) were re-sized (narrowed) to fit the new text.
#
# This script walks every worksheet in the workbook, replaces the old
# status text wherever it is found, and narrows the column(s) that held
# that text.

$wb = $excel.ActiveWorkbook

$oldStatusText = "Ready for handoff"
$newStatusText = "In Translation"

# Target width taken from the updated report (character-width units, as
# used by Range.ColumnWidth).
$newColumnWidth = 13.4101845877511

# Columns get resized only if they currently hold the status text's old
# width (~17.22 character-width units / ColumnWidth ~16.33). Use a small
# tolerance since ColumnWidth is a floating point, quantized value.
$oldColumnWidthApprox = 16.33
$widthTolerance = 0.02

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange

    # --- 1. Replace the status text -------------------------------------
    foreach ($cell in $used.Cells) {
        # Cast through [string] explicitly: for boolean-valued cells,
        # comparing the raw (non-string) Text/Value2 against a string
        # literal can coerce the string to a boolean instead of coercing
        # the boolean to a string, giving false positives.
        $cellText = [string]$cell.Text
        if ($cellText -eq $oldStatusText) {
            $cell.Value2 = $newStatusText
        }
    }

    # --- 2. Narrow the column(s) that hold the status text --------------
    foreach ($col in $used.Columns) {
        $width = $col.ColumnWidth
        if ([math]::Abs($width - $oldColumnWidthApprox) -lt $widthTolerance) {
            $col.ColumnWidth = $newColumnWidth
        }
    }
}
